# Update Lcn2-Lrp2 LR-pair sheet with new TPM-based NATMI output.
# - Adds a new "Resolving-Mac" sending-cluster row (row 7).
# - Re-points existing rows' Ligand/Receptor symbol cells (B:C), which shift
#   shared-string indices once "Resolving-Mac" is introduced into the pool.
# - Refreshes the recomputed statistics columns (E:T) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T)
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
}

Set-Row 2 "ECs" "Lcn2" "Lrp2" "MuSCs" `
    2 0.6666666666666666 `
    0.195866 0.5875980000000001 `
    0.00204520224202265 0.00204520224202265 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    0.021550744248 0.193956698232 `
    0.00204520224202265 0.00204520224202265

Set-Row 3 "FAPs" "Lcn2" "Lrp2" "MuSCs" `
    3 1 `
    2.208069666666666 6.624209 `
    0.02305631928363714 0.02305631928363714 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    0.242949489284 2.186545403556 `
    0.02305631928363714 0.02305631928363714

Set-Row 4 "Inflammatory-Mac" "Lcn2" "Lrp2" "MuSCs" `
    3 1 `
    8.31967 24.95901 `
    0.08687269733842821 0.08687269733842821 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    0.9153966507600001 8.23856985684 `
    0.08687269733842821 0.08687269733842821

Set-Row 5 "MuSCs" "Lcn2" "Lrp2" "MuSCs" `
    2 0.6666666666666666 `
    0.1125653333333333 0.337696 `
    0.001175389664910501 0.001175389664910501 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    0.012385338496 0.111468046464 `
    0.001175389664910501 0.001175389664910501

Set-Row 6 "Neutrophils" "Lcn2" "Lrp2" "MuSCs" `
    3 1 `
    84.705523 254.116569 `
    0.8844818679673917 0.8844818679673917 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    9.319979284644001 83.879813561796 `
    0.8844818679673917 0.8844818679673917

Set-Row 7 "Resolving-Mac" "Lcn2" "Lrp2" "MuSCs" `
    2 0.6666666666666666 `
    0.22683 0.68049 `
    0.002368523503609599 0.002368523503609599 `
    3 1 `
    0.110028 0.330084 `
    1 1 `
    0.02495765124 0.22461886116 `
    0.002368523503609599 0.002368523503609599
